$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
"42+31=",
"15-5=",
"89-64=",
"58-24=",
"24+52=",
"61-0=",
"40-38=",
"10+86=",
"39+37=",
"73+0=",
"95-16=",
"5+41=",
"78-18=",
"35+63=",
"89-25=",
"98-7=",
"24-8=",
"64-2=",
"13+73=",
"80-28=",
"84-16=",
"0+53=",
"29-9=",
"29+60=",
"10+88=",
"17+36=",
"60+1=",
"94-11=",
"80+15=",
"40-14=",
"44+47=",
"62-28=",
"71+17=",
"11+52=",
"79-33=",
"24+56=",
"34-18=",
"66-48=",
"6+21=",
"95-32=",
"28-13=",
"18-6=",
"30+67=",
"65-46=",
"87-87=",
"73-44=",
"8+81=",
"99-59=",
"3+41=",
"72-42=",
"95-13=",
"61-7=",
"74-33=",
"91-34=",
"46-27=",
"76-47=",
"21+39=",
"86-48=",
"51+20=",
"96-72=",
"62+23=",
"74-24=",
"62-38=",
"27+37=",
"12+74=",
"37-20=",
"89-82=",
"43-28=",
"72+11=",
"70-27=",
"22-13=",
"89-19=",
"47-32=",
"20+35=",
"47-17=",
"30-21=",
"37-14=",
"92-76=",
"58+3=",
"96-42=",
"51+40=",
"8+35=",
"43+42=",
"83-51=",
"6+58=",
"76-41=",
"92+2=",
"75-12=",
"78+9=",
"49-14=",
"70+2=",
"11+48=",
"17+71=",
"67-46=",
"61-22=",
"36+55=",
"56-53=",
"7+8=",
"68-36=",
"7-5="
)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
if ($rows * $cols -ne $values.Count) {
    throw "Expected $($values.Count) cells but table has $rows x $cols = $($rows * $cols)"
}
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$i]
        $i++
    }
}
Write-Output "done: $i cells updated"